$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.009228348731995
$ws.Range("B1").Value = 3.477327585220337
$ws.Range("C1").Value = 3.929091215133667
$ws.Range("D1").Value = 3.136507749557495
$ws.Range("E1").Value = 1.304911494255066
